# Edit script: Add ABW and AGW experimental data rows, update calibration
# counts/errors derived from the new data, and update sheet selections to
# match the final authored state. Downstream formulas (Parameters!B8:B9,
# Bottle Results, Averaged Results, Count->Actual Activity!F:G) recalc
# automatically from the underlying data changes made below.

$wb = $excel.ActiveWorkbook

$wsParams = $wb.Worksheets.Item("Parameters")
$wsScint  = $wb.Worksheets.Item("Scintillation Counter Results")
$wsCal    = $wb.Worksheets.Item("Calibration Data")
$wsCount  = $wb.Worksheets.Item("Count->Actual Activity")
$wsBottle = $wb.Worksheets.Item("Bottle Results")
$wsAvg    = $wb.Worksheets.Item("Averaged Results")

# ---------------------------------------------------------------------
# 1. Append 75 new raw scintillation-counter rows (rows 32-106) to the
#    "Scintillation Counter Results" sheet with the newly measured AGW /
#    ABW / FHY counting data.
# ---------------------------------------------------------------------
$wsScint.Range("A32").Value = 42926.624305555553
$wsScint.Range("B32").Value = "AGW Initial Stock"
$wsScint.Range("C32").Value = 2430.6999999999998
$wsScint.Range("D32").Value = 1.28
$wsScint.Range("E32").Value = 0
$wsScint.Range("F32").Value = 43
$wsScint.Range("A33").Value = 42926.624305555553
$wsScint.Range("B33").Value = "RaGlassAGW_1A"
$wsScint.Range("C33").Value = 2308
$wsScint.Range("D33").Value = 1.32
$wsScint.Range("E33").Value = 0
$wsScint.Range("F33").Value = 53.61
$wsScint.Range("A34").Value = 42926.624305555553
$wsScint.Range("B34").Value = "RaGlassAGW_1A (5mL)"
$wsScint.Range("C34").Value = 1201.7
$wsScint.Range("D34").Value = 1.82
$wsScint.Range("E34").Value = 0.01
$wsScint.Range("F34").Value = 64.239999999999995
$wsScint.Range("A35").Value = 42926.624305555553
$wsScint.Range("B35").Value = "RaGlassAGW_1B"
$wsScint.Range("C35").Value = 2231.4
$wsScint.Range("D35").Value = 1.34
$wsScint.Range("E35").Value = 0
$wsScint.Range("F35").Value = 74.86
$wsScint.Range("A36").Value = 42926.624305555553
$wsScint.Range("B36").Value = "RaGlassAGW_1C"
$wsScint.Range("C36").Value = 2300.8000000000002
$wsScint.Range("D36").Value = 1.32
$wsScint.Range("E36").Value = 0
$wsScint.Range("F36").Value = 85.49
$wsScint.Range("A37").Value = 42926.624305555553
$wsScint.Range("B37").Value = "RaMontAGW_1A"
$wsScint.Range("C37").Value = 1692.8
$wsScint.Range("D37").Value = 1.54
$wsScint.Range("E37").Value = 0
$wsScint.Range("F37").Value = 96.11
$wsScint.Range("A38").Value = 42926.624305555553
$wsScint.Range("B38").Value = "RaMontAGW_1B"
$wsScint.Range("C38").Value = 1725.1
$wsScint.Range("D38").Value = 1.52
$wsScint.Range("E38").Value = 0.01
$wsScint.Range("F38").Value = 106.74
$wsScint.Range("A39").Value = 42926.624305555553
$wsScint.Range("B39").Value = "RaMontAGW_1C"
$wsScint.Range("C39").Value = 1693.9
$wsScint.Range("D39").Value = 1.54
$wsScint.Range("E39").Value = 0.01
$wsScint.Range("F39").Value = 117.36
$wsScint.Range("A40").Value = 42926.624305555553
$wsScint.Range("B40").Value = "AGW Final Stock"
$wsScint.Range("C40").Value = 2363.9
$wsScint.Range("D40").Value = 1.3
$wsScint.Range("E40").Value = 0
$wsScint.Range("F40").Value = 128
$wsScint.Range("A41").Value = 42926.624305555553
$wsScint.Range("B41").Value = "RaFHYAGW_1A"
$wsScint.Range("C41").Value = 1890.2
$wsScint.Range("D41").Value = 1.45
$wsScint.Range("E41").Value = 0.01
$wsScint.Range("F41").Value = 138.72999999999999
$wsScint.Range("A42").Value = 42926.624305555553
$wsScint.Range("B42").Value = "RaFHYAGW_1B"
$wsScint.Range("C42").Value = 1942.3
$wsScint.Range("D42").Value = 1.44
$wsScint.Range("E42").Value = 0
$wsScint.Range("F42").Value = 149.36000000000001
$wsScint.Range("A43").Value = 42926.624305555553
$wsScint.Range("B43").Value = "RaFHYAGW_1C"
$wsScint.Range("C43").Value = 1916.8
$wsScint.Range("D43").Value = 1.44
$wsScint.Range("E43").Value = 0
$wsScint.Range("F43").Value = 159.97
$wsScint.Range("A44").Value = 42926.624305555553
$wsScint.Range("B44").Value = "RaGOEAGW_1A"
$wsScint.Range("C44").Value = 2093
$wsScint.Range("D44").Value = 1.38
$wsScint.Range("E44").Value = 0
$wsScint.Range("F44").Value = 170.61
$wsScint.Range("A45").Value = 42926.624305555553
$wsScint.Range("B45").Value = "RaGOEAGW_1B"
$wsScint.Range("C45").Value = 2118.1
$wsScint.Range("D45").Value = 1.37
$wsScint.Range("E45").Value = 0
$wsScint.Range("F45").Value = 181.24
$wsScint.Range("A46").Value = 42926.624305555553
$wsScint.Range("B46").Value = "RaGOWAGW_1C"
$wsScint.Range("C46").Value = 2120.1999999999998
$wsScint.Range("D46").Value = 1.37
$wsScint.Range("E46").Value = 0.01
$wsScint.Range("F46").Value = 191.88
$wsScint.Range("A47").Value = 42923.320833333331
$wsScint.Range("B47").Value = "AGW Initial Stock"
$wsScint.Range("C47").Value = 2449.4
$wsScint.Range("D47").Value = 1.28
$wsScint.Range("E47").Value = 0
$wsScint.Range("F47").Value = 42.93
$wsScint.Range("A48").Value = 42923.320833333331
$wsScint.Range("B48").Value = "RaGlassAGW_1A"
$wsScint.Range("C48").Value = 2326.6999999999998
$wsScint.Range("D48").Value = 1.31
$wsScint.Range("E48").Value = 0
$wsScint.Range("F48").Value = 53.56
$wsScint.Range("A49").Value = 42923.320833333331
$wsScint.Range("B49").Value = "RaGlassAGW_1A (5mL)"
$wsScint.Range("C49").Value = 1202.2
$wsScint.Range("D49").Value = 1.82
$wsScint.Range("E49").Value = 0.01
$wsScint.Range("F49").Value = 64.180000000000007
$wsScint.Range("A50").Value = 42923.320833333331
$wsScint.Range("B50").Value = "RaGlassAGW_1B"
$wsScint.Range("C50").Value = 2256
$wsScint.Range("D50").Value = 1.33
$wsScint.Range("E50").Value = 0
$wsScint.Range("F50").Value = 74.81
$wsScint.Range("A51").Value = 42923.320833333331
$wsScint.Range("B51").Value = "RaGlassAGW_1C"
$wsScint.Range("C51").Value = 2330.3000000000002
$wsScint.Range("D51").Value = 1.31
$wsScint.Range("E51").Value = 0
$wsScint.Range("F51").Value = 85.43
$wsScint.Range("A52").Value = 42923.320833333331
$wsScint.Range("B52").Value = "RaMontAGW_1A"
$wsScint.Range("C52").Value = 1715.7
$wsScint.Range("D52").Value = 1.53
$wsScint.Range("E52").Value = 0.01
$wsScint.Range("F52").Value = 96.06
$wsScint.Range("A53").Value = 42923.320833333331
$wsScint.Range("B53").Value = "RaMontAGW_1B"
$wsScint.Range("C53").Value = 1705
$wsScint.Range("D53").Value = 1.53
$wsScint.Range("E53").Value = 0.01
$wsScint.Range("F53").Value = 106.68
$wsScint.Range("A54").Value = 42923.320833333331
$wsScint.Range("B54").Value = "RaMontAGW_1C"
$wsScint.Range("C54").Value = 1700.7
$wsScint.Range("D54").Value = 1.53
$wsScint.Range("E54").Value = 0.01
$wsScint.Range("F54").Value = 117.31
$wsScint.Range("A55").Value = 42923.320833333331
$wsScint.Range("B55").Value = "AGW Final Stock"
$wsScint.Range("C55").Value = 2407.1
$wsScint.Range("D55").Value = 1.29
$wsScint.Range("E55").Value = 0
$wsScint.Range("F55").Value = 127.95
$wsScint.Range("A56").Value = 42923.320833333331
$wsScint.Range("B56").Value = "RaFHYAGW_1A"
$wsScint.Range("C56").Value = 1907.9
$wsScint.Range("D56").Value = 1.45
$wsScint.Range("E56").Value = 0.02
$wsScint.Range("F56").Value = 138.69
$wsScint.Range("A57").Value = 42923.320833333331
$wsScint.Range("B57").Value = "RaFHYAGW_1B"
$wsScint.Range("C57").Value = 1962.1
$wsScint.Range("D57").Value = 1.43
$wsScint.Range("E57").Value = 0
$wsScint.Range("F57").Value = 149.32
$wsScint.Range("A58").Value = 42923.320833333331
$wsScint.Range("B58").Value = "RaFHYAGW_1C"
$wsScint.Range("C58").Value = 1957.3
$wsScint.Range("D58").Value = 1.43
$wsScint.Range("E58").Value = 0.01
$wsScint.Range("F58").Value = 159.94
$wsScint.Range("A59").Value = 42923.320833333331
$wsScint.Range("B59").Value = "RaGOEAGW_1A"
$wsScint.Range("C59").Value = 2130.9
$wsScint.Range("D59").Value = 1.37
$wsScint.Range("E59").Value = 0
$wsScint.Range("F59").Value = 170.58
$wsScint.Range("A60").Value = 42923.320833333331
$wsScint.Range("B60").Value = "RaGOEAGW_1B"
$wsScint.Range("C60").Value = 2153.8000000000002
$wsScint.Range("D60").Value = 1.38
$wsScint.Range("E60").Value = 0
$wsScint.Range("F60").Value = 181.21
$wsScint.Range("A61").Value = 42923.320833333331
$wsScint.Range("B61").Value = "RaGOWAGW_1C"
$wsScint.Range("C61").Value = 2171.5
$wsScint.Range("D61").Value = 1.36
$wsScint.Range("E61").Value = 0.01
$wsScint.Range("F61").Value = 191.85
$wsScint.Range("A62").Value = 42927.490277777775
$wsScint.Range("B62").Value = "AGW Initial Stock"
$wsScint.Range("C62").Value = 2420.1999999999998
$wsScint.Range("D62").Value = 1.29
$wsScint.Range("E62").Value = 0
$wsScint.Range("F62").Value = 42.9
$wsScint.Range("A63").Value = 42927.490277777775
$wsScint.Range("B63").Value = "RaGlassAGW_1A"
$wsScint.Range("C63").Value = 2303.6
$wsScint.Range("D63").Value = 1.32
$wsScint.Range("E63").Value = 0
$wsScint.Range("F63").Value = 53.51
$wsScint.Range("A64").Value = 42927.490277777775
$wsScint.Range("B64").Value = "RaGlassAGW_1A (5mL)"
$wsScint.Range("C64").Value = 1219.9000000000001
$wsScint.Range("D64").Value = 1.81
$wsScint.Range("E64").Value = 0.01
$wsScint.Range("F64").Value = 64.150000000000006
$wsScint.Range("A65").Value = 42927.490277777775
$wsScint.Range("B65").Value = "RaGlassAGW_1B"
$wsScint.Range("C65").Value = 2257.9
$wsScint.Range("D65").Value = 1.33
$wsScint.Range("E65").Value = 0
$wsScint.Range("F65").Value = 74.760000000000005
$wsScint.Range("A66").Value = 42927.490277777775
$wsScint.Range("B66").Value = "RaGlassAGW_1C"
$wsScint.Range("C66").Value = 2303.5
$wsScint.Range("D66").Value = 1.32
$wsScint.Range("E66").Value = 0
$wsScint.Range("F66").Value = 85.39
$wsScint.Range("A67").Value = 42927.490277777775
$wsScint.Range("B67").Value = "RaMontAGW_1A"
$wsScint.Range("C67").Value = 1709.1
$wsScint.Range("D67").Value = 1.53
$wsScint.Range("E67").Value = 0.01
$wsScint.Range("F67").Value = 96.01
$wsScint.Range("A68").Value = 42927.490277777775
$wsScint.Range("B68").Value = "RaMontAGW_1B"
$wsScint.Range("C68").Value = 1716.8
$wsScint.Range("D68").Value = 1.53
$wsScint.Range("E68").Value = 0.01
$wsScint.Range("F68").Value = 106.64
$wsScint.Range("A69").Value = 42927.490277777775
$wsScint.Range("B69").Value = "RaMontAGW_1C"
$wsScint.Range("C69").Value = 1708.2
$wsScint.Range("D69").Value = 1.53
$wsScint.Range("E69").Value = 0.01
$wsScint.Range("F69").Value = 117.26
$wsScint.Range("A70").Value = 42927.490277777775
$wsScint.Range("B70").Value = "AGW Final Stock"
$wsScint.Range("C70").Value = 2416.5
$wsScint.Range("D70").Value = 1.29
$wsScint.Range("E70").Value = 0
$wsScint.Range("F70").Value = 127.9
$wsScint.Range("A71").Value = 42927.490277777775
$wsScint.Range("B71").Value = "RaFHYAGW_1A"
$wsScint.Range("C71").Value = 1915.8
$wsScint.Range("D71").Value = 1.44
$wsScint.Range("E71").Value = 0.01
$wsScint.Range("F71").Value = 138.63
$wsScint.Range("A72").Value = 42927.490277777775
$wsScint.Range("B72").Value = "RaFHYAGW_1B"
$wsScint.Range("C72").Value = 1914.5
$wsScint.Range("D72").Value = 1.45
$wsScint.Range("E72").Value = 0
$wsScint.Range("F72").Value = 149.26
$wsScint.Range("A73").Value = 42927.490277777775
$wsScint.Range("B73").Value = "RaFHYAGW_1C"
$wsScint.Range("C73").Value = 1939.7
$wsScint.Range("D73").Value = 1.44
$wsScint.Range("E73").Value = 0
$wsScint.Range("F73").Value = 159.88999999999999
$wsScint.Range("A74").Value = 42927.490277777775
$wsScint.Range("B74").Value = "RaGOEAGW_1A"
$wsScint.Range("C74").Value = 2135
$wsScint.Range("D74").Value = 1.37
$wsScint.Range("E74").Value = 0
$wsScint.Range("F74").Value = 170.53
$wsScint.Range("A75").Value = 42927.490277777775
$wsScint.Range("B75").Value = "RaGOEAGW_1B"
$wsScint.Range("C75").Value = 2106
$wsScint.Range("D75").Value = 1.38
$wsScint.Range("E75").Value = 0
$wsScint.Range("F75").Value = 181.16
$wsScint.Range("A76").Value = 42927.490277777775
$wsScint.Range("B76").Value = "RaGOWAGW_1C"
$wsScint.Range("C76").Value = 2144.8000000000002
$wsScint.Range("D76").Value = 1.37
$wsScint.Range("E76").Value = 0
$wsScint.Range("F76").Value = 191.8
$wsScint.Range("A77").Value = 42928.348611111112
$wsScint.Range("B77").Value = "AGW Initial Stock"
$wsScint.Range("C77").Value = 2371.6
$wsScint.Range("D77").Value = 1.3
$wsScint.Range("E77").Value = 0
$wsScint.Range("F77").Value = 42.9
$wsScint.Range("A78").Value = 42928.348611111112
$wsScint.Range("B78").Value = "RaGlassAGW_1A"
$wsScint.Range("C78").Value = 2252.1
$wsScint.Range("D78").Value = 1.33
$wsScint.Range("E78").Value = 0
$wsScint.Range("F78").Value = 53.51
$wsScint.Range("A79").Value = 42928.348611111112
$wsScint.Range("B79").Value = "RaGlassAGW_1A (5mL)"
$wsScint.Range("C79").Value = 1173.3
$wsScint.Range("D79").Value = 1.85
$wsScint.Range("E79").Value = 0.01
$wsScint.Range("F79").Value = 64.14
$wsScint.Range("A80").Value = 42928.348611111112
$wsScint.Range("B80").Value = "RaGlassAGW_1B"
$wsScint.Range("C80").Value = 2265
$wsScint.Range("D80").Value = 1.33
$wsScint.Range("E80").Value = 0
$wsScint.Range("F80").Value = 74.760000000000005
$wsScint.Range("A81").Value = 42928.348611111112
$wsScint.Range("B81").Value = "RaGlassAGW_1C"
$wsScint.Range("C81").Value = 2266
$wsScint.Range("D81").Value = 1.33
$wsScint.Range("E81").Value = 0
$wsScint.Range("F81").Value = 85.39
$wsScint.Range("A82").Value = 42928.348611111112
$wsScint.Range("B82").Value = "RaMontAGW_1A"
$wsScint.Range("C82").Value = 1676.7
$wsScint.Range("D82").Value = 1.54
$wsScint.Range("E82").Value = 0
$wsScint.Range("F82").Value = 96.01
$wsScint.Range("A83").Value = 42928.348611111112
$wsScint.Range("B83").Value = "RaMontAGW_1B"
$wsScint.Range("C83").Value = 1706.8
$wsScint.Range("D83").Value = 1.53
$wsScint.Range("E83").Value = 0.01
$wsScint.Range("F83").Value = 106.63
$wsScint.Range("A84").Value = 42928.348611111112
$wsScint.Range("B84").Value = "RaMontAGW_1C"
$wsScint.Range("C84").Value = 1669.4
$wsScint.Range("D84").Value = 1.55
$wsScint.Range("E84").Value = 0.01
$wsScint.Range("F84").Value = 117.26
$wsScint.Range("A85").Value = 42928.348611111112
$wsScint.Range("B85").Value = "AGW Final Stock"
$wsScint.Range("C85").Value = 2390.1
$wsScint.Range("D85").Value = 1.29
$wsScint.Range("E85").Value = 0
$wsScint.Range("F85").Value = 127.9
$wsScint.Range("A86").Value = 42928.348611111112
$wsScint.Range("B86").Value = "RaFHYAGW_1A"
$wsScint.Range("C86").Value = 1871.6
$wsScint.Range("D86").Value = 1.46
$wsScint.Range("E86").Value = 0.01
$wsScint.Range("F86").Value = 138.63
$wsScint.Range("A87").Value = 42928.348611111112
$wsScint.Range("B87").Value = "RaFHYAGW_1B"
$wsScint.Range("C87").Value = 1888.8
$wsScint.Range("D87").Value = 1.46
$wsScint.Range("E87").Value = 0
$wsScint.Range("F87").Value = 149.24
$wsScint.Range("A88").Value = 42928.348611111112
$wsScint.Range("B88").Value = "RaFHYAGW_1C"
$wsScint.Range("C88").Value = 1928.8
$wsScint.Range("D88").Value = 1.44
$wsScint.Range("E88").Value = 0
$wsScint.Range("F88").Value = 159.88
$wsScint.Range("A89").Value = 42928.348611111112
$wsScint.Range("B89").Value = "RaGOEAGW_1A"
$wsScint.Range("C89").Value = 2091.6999999999998
$wsScint.Range("D89").Value = 1.38
$wsScint.Range("E89").Value = 0
$wsScint.Range("F89").Value = 170.51
$wsScint.Range("A90").Value = 42928.348611111112
$wsScint.Range("B90").Value = "RaGOEAGW_1B"
$wsScint.Range("C90").Value = 2119.3000000000002
$wsScint.Range("D90").Value = 1.37
$wsScint.Range("E90").Value = 0
$wsScint.Range("F90").Value = 181.14
$wsScint.Range("A91").Value = 42928.348611111112
$wsScint.Range("B91").Value = "RaGOWAGW_1C"
$wsScint.Range("C91").Value = 2127
$wsScint.Range("D91").Value = 1.37
$wsScint.Range("E91").Value = 0
$wsScint.Range("F91").Value = 191.78
$wsScint.Range("A92").Value = 42926.624305555553
$wsScint.Range("B92").Value = "AGW Initial Stock"
$wsScint.Range("C92").Value = 2430.6999999999998
$wsScint.Range("D92").Value = 1.28
$wsScint.Range("E92").Value = 0
$wsScint.Range("F92").Value = 43
$wsScint.Range("A93").Value = 42926.624305555553
$wsScint.Range("B93").Value = "RaGlassAGW_1A"
$wsScint.Range("C93").Value = 2308.6999999999998
$wsScint.Range("D93").Value = 1.32
$wsScint.Range("E93").Value = 0
$wsScint.Range("F93").Value = 53.61
$wsScint.Range("A94").Value = 42926.624305555553
$wsScint.Range("B94").Value = "RaGlassAGW_1A (5mL)"
$wsScint.Range("C94").Value = 1201.7
$wsScint.Range("D94").Value = 1.82
$wsScint.Range("E94").Value = 0.01
$wsScint.Range("F94").Value = 64.239999999999995
$wsScint.Range("A95").Value = 42926.624305555553
$wsScint.Range("B95").Value = "RaGlassAGW_1B"
$wsScint.Range("C95").Value = 2231.4
$wsScint.Range("D95").Value = 1.34
$wsScint.Range("E95").Value = 0
$wsScint.Range("F95").Value = 74.86
$wsScint.Range("A96").Value = 42926.624305555553
$wsScint.Range("B96").Value = "RaGlassAGW_1C"
$wsScint.Range("C96").Value = 2300.8000000000002
$wsScint.Range("D96").Value = 1.32
$wsScint.Range("E96").Value = 0
$wsScint.Range("F96").Value = 85.49
$wsScint.Range("A97").Value = 42926.624305555553
$wsScint.Range("B97").Value = "RaMontAGW_1A"
$wsScint.Range("C97").Value = 1692.8
$wsScint.Range("D97").Value = 1.54
$wsScint.Range("E97").Value = 0
$wsScint.Range("F97").Value = 96.11
$wsScint.Range("A98").Value = 42926.624305555553
$wsScint.Range("B98").Value = "RaMontAGW_1B"
$wsScint.Range("C98").Value = 1725.1
$wsScint.Range("D98").Value = 1.52
$wsScint.Range("E98").Value = 0.01
$wsScint.Range("F98").Value = 106.74
$wsScint.Range("A99").Value = 42926.624305555553
$wsScint.Range("B99").Value = "RaMontAGW_1C"
$wsScint.Range("C99").Value = 1693.9
$wsScint.Range("D99").Value = 1.54
$wsScint.Range("E99").Value = 0.01
$wsScint.Range("F99").Value = 117.36
$wsScint.Range("A100").Value = 42926.624305555553
$wsScint.Range("B100").Value = "AGW Final Stock"
$wsScint.Range("C100").Value = 2363.9
$wsScint.Range("D100").Value = 1.3
$wsScint.Range("E100").Value = 0
$wsScint.Range("F100").Value = 128
$wsScint.Range("A101").Value = 42926.624305555553
$wsScint.Range("B101").Value = "RaFHYAGW_1A"
$wsScint.Range("C101").Value = 1890.2
$wsScint.Range("D101").Value = 1.45
$wsScint.Range("E101").Value = 0.01
$wsScint.Range("F101").Value = 138.72999999999999
$wsScint.Range("A102").Value = 42926.624305555553
$wsScint.Range("B102").Value = "RaFHYAGW_1B"
$wsScint.Range("C102").Value = 1942.3
$wsScint.Range("D102").Value = 1.44
$wsScint.Range("E102").Value = 0
$wsScint.Range("F102").Value = 149.36000000000001
$wsScint.Range("A103").Value = 42926.624305555553
$wsScint.Range("B103").Value = "RaFHYAGW_1C"
$wsScint.Range("C103").Value = 1916.8
$wsScint.Range("D103").Value = 1.44
$wsScint.Range("E103").Value = 0
$wsScint.Range("F103").Value = 159.97
$wsScint.Range("A104").Value = 42926.624305555553
$wsScint.Range("B104").Value = "RaGOEAGW_1A"
$wsScint.Range("C104").Value = 2093
$wsScint.Range("D104").Value = 1.38
$wsScint.Range("E104").Value = 0
$wsScint.Range("F104").Value = 170.61
$wsScint.Range("A105").Value = 42926.624305555553
$wsScint.Range("B105").Value = "RaGOEAGW_1B"
$wsScint.Range("C105").Value = 2118.1
$wsScint.Range("D105").Value = 1.37
$wsScint.Range("E105").Value = 0
$wsScint.Range("F105").Value = 181.24
$wsScint.Range("A106").Value = 42926.624305555553
$wsScint.Range("B106").Value = "RaGOWAGW_1C"
$wsScint.Range("C106").Value = 2120.1999999999998
$wsScint.Range("D106").Value = 1.37
$wsScint.Range("E106").Value = 0.01
$wsScint.Range("F106").Value = 191.88

# ---------------------------------------------------------------------
# 2. Update the "Count->Actual Activity" mean-CPS counts (column C) and
#    their propagated errors (column D) to reflect the larger pool of
#    raw measurements now available for each sample on the
#    "Scintillation Counter Results" sheet. Downstream activity
#    calculations (columns F/G), Parameters!B8:B9, Bottle Results, and
#    Averaged Results all recalculate automatically from these values.
# ---------------------------------------------------------------------
$ws4 = $wsCount
$ws4.Range("C2").Value = 40.454047619047607
$ws4.Range("D2").Value = 0.51896763945578228
$ws4.Range("C3").Value = 38.326428571428558
$ws4.Range("D3").Value = 0.50536133673469386
$ws4.Range("C4").Value = 19.988571428571429
$ws4.Range("D4").Value = 0.36464865306122451
$ws4.Range("C5").Value = 37.770238095238092
$ws4.Range("D5").Value = 0.5023441666666667
$ws4.Range("C6").Value = 38.499523809523808
$ws4.Range("D6").Value = 0.50709372789115648
$ws4.Range("C7").Value = 28.39595238095238
$ws4.Range("D7").Value = 0.43526938435374152
$ws4.Range("C8").Value = 28.653571428571428
$ws4.Range("D8").Value = 0.43676229591836752
$ws4.Range("C9").Value = 28.26428571428572
$ws4.Range("D9").Value = 0.43405867346938781
$ws4.Range("C10").Value = 39.940714285714293
$ws4.Range("D10").Value = 0.51580579591836739
$ws4.Range("C11").Value = 31.722380952380949
$ws4.Range("D11").Value = 0.45952134693877539
$ws4.Range("C12").Value = 32.367142857142859
$ws4.Range("D12").Value = 0.46562446938775509
$ws4.Range("C13").Value = 32.304285714285712
$ws4.Range("D13").Value = 0.4637972448979592
$ws4.Range("C14").Value = 35.265000000000001
$ws4.Range("D14").Value = 0.48464185714285718
$ws4.Range("C15").Value = 35.556904761904768
$ws4.Range("D15").Value = 0.48712959523809529
$ws4.Range("C16").Value = 35.808809523809529
$ws4.Range("D16").Value = 0.48853447278911571

# ---------------------------------------------------------------------
# 3. Restore / update the selection shown on each worksheet to match the
#    final authored state, finishing on "Averaged Results" so that sheet
#    remains the active tab (as in the original workbook).
# ---------------------------------------------------------------------
$wsParams.Range("D8:E9").Select()

$wsScint.Range("C107").Select()

$wsCount.Range("C2:D16").Select()

$wsBottle.Range("D8").Select()

$wsAvg.Range("I11").Select()
